$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("length_of_stay_categories")

# --- Insert the two new columns ---------------------------------------
# New "length_of_stay_simple_three_days" column goes in before the old B
# (two_weeks), pushing two_weeks/week/three right by one.
$ws.Columns("B:B").Insert()

# New "length_of_stay_simple_three_days_order" column goes in before the
# old G (three_order), which by now lives in column H after the first
# insert; pushing it right by one to I.
$ws.Columns("H:H").Insert()

# --- Formatting ----------------------------------------------------------
# Match the text number format ("@") used by the other string columns so
# the new cells pick up the same cell style (s="3") as their neighbours,
# and so that the numeric-looking labels in column H ("1"/"2") are stored
# as text (shared strings) instead of being auto-coerced to numbers.
$ws.Range("B1:B60").NumberFormat = "@"
$ws.Range("H1:H60").NumberFormat = "@"

# --- Header row (row 1) ----------------------------------------------
$ws.Range("B1").Value = "length_of_stay_simple_three_days"
$ws.Range("H1").Value = "length_of_stay_simple_three_days_order"

# --- Body values -----------------------------------------------------
# length_of_stay (A) 1-6  -> three_days "1-3", three_days_order "1"
# length_of_stay (A) 7-59 -> three_days "4+",  three_days_order "2"
$ws.Range("B2:B7").Value = "1-3"
$ws.Range("H2:H7").Value = "1"

$ws.Range("B8:B60").Value = "4+"
$ws.Range("H8:H60").Value = "2"

# Column widths: B-D (20.5 chars) share the width of the former B-C block,
# E keeps 17.5, and the new H also gets 20.5. The COM ColumnWidth setter
# adds a fixed ~0.8333 character padding versus the stored OOXML width, so
# we compensate to land exactly on 20.5 / 17.5 after round-tripping.
$ws.Columns("B:D").ColumnWidth = 19.666666666666668
$ws.Columns("E:E").ColumnWidth = 16.666666666666668
$ws.Columns("H:H").ColumnWidth = 19.666666666666668

# --- Selection ----------------------------------------------------------
$ws.Range("B5:B60").Select()
